# This script updates the cryptos worksheet (rows 2-51, columns B-E) to the
# new snapshot values. Column D (Price) holds values that look numeric (e.g.
# "573.91") but must remain plain text, matching the workbook's original
# inlineStr encoding -- Excel's COM layer auto-coerces such strings to real
# numbers on assignment, so we format the column as Text first, write the
# values, then restore the cell style so no stray numeric formatting lingers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "62.637.65", "  +4.06%  ")
    ,@("Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "2.413.37", "  +1.17%  ")
    ,@("TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.00", "  +0.20%  ")
    ,@("BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "573.91", "  +2.40%  ")
    ,@("Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "145.83", "  +5.55%  ")
    ,@("USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "0.997", "  -0.39%  ")
    ,@("XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.540", "  +0.58%  ")
    ,@("LidoStakedEther", "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth", "2.439.50", "  +2.20%  ")
    ,@("Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.111", "  +5.37%  ")
    ,@("TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.161", "  +0.85%  ")
    ,@("Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "5.25", "  +3.71%  ")
    ,@("Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.352", "  +4.54%  ")
    ,@("Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "26.88", "  +4.94%  ")
    ,@("ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.0000179", "  +8.42%  ")
    ,@("WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "2.877.21", "  +2.13%  ")
    ,@("WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "62.339.60", "  +3.51%  ")
    ,@("WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "2.440.79", "  +1.91%  ")
    ,@("Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "7.97", "  -5.07%  ")
    ,@("Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "10.91", "  +3.33%  ")
    ,@("BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "326.67", "  +0.73%  ")
    ,@("Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "4.13", "  +2.78%  ")
    ,@("SuiNetwork", "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui", "2.04", "  +14.70%  ")
    ,@("Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "1.00", "  -0.06%  ")
    ,@("Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "65.63", "  +1.79%  ")
    ,@("Bittensor", "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao", "618.94", "  +11.94%  ")
    ,@("Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "8.38", "  +5.34%  ")
    ,@("PEPE", "https://coinranking.com/coin/03WI8NQPF+pepe-pepe", "0.0₃0990", "  +10.09%  ")
    ,@("WrappedeETH", "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth", "2.534.64", "  +1.04%  ")
    ,@("InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "8.12", "  +2.61%  ")
    ,@("Binance-PegBSC-USD", "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd", "0.952", "  -4.88%  ")
    ,@("Fetch.AI", "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet", "1.41", "  +9.66%  ")
    ,@("Kaspa", "https://coinranking.com/coin/V8GxkwWow+kaspa-kas", "0.139", "  +7.01%  ")
    ,@("PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "1.84", "  +2.68%  ")
    ,@("ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "1.48", "  +5.40%  ")
    ,@("FirstDigitalUSD", "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd", "0.994", "  -0.50%  ")
    ,@("NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", "4.75", "  +5.82%  ")
    ,@("PolygonEcosystemToken", "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol", "0.373", "  +1.86%  ")
    ,@("RenderToken", "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render", "5.40", "  +7.86%  ")
    ,@("Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "151.28", "  -1.36%  ")
    ,@("EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "18.58", "  +1.97%  ")
    ,@("dogwifhat", "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif", "2.76", "  +21.30%  ")
    ,@("Stacks", "https://coinranking.com/coin/mMPrMcB7+stacks-stx", "1.76", "  +7.90%  ")
    ,@("OKB", "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb", "42.34", "  +3.03%  ")
    ,@("USDe", "https://coinranking.com/coin/exbfr2U-0+usde-usde", "0.999", "  +0.01%  ")
    ,@("BabyDogeCoin", "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge", "0.0₆0282", "  +1.08%  ")
    ,@("Aave", "https://coinranking.com/coin/ixgUfzmLR+aave-aave", "144.54", "  +1.55%  ")
    ,@("Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "3.58", "  +2.81%  ")
    ,@("InjectiveProtocol", "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj", "20.28", "  +7.87%  ")
    ,@("Mantle", "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt", "0.601", "  +2.42%  ")
    ,@("Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.0515", "  +3.40%  ")
)

$firstRow = 2
$lastRow = $firstRow + $data.Count - 1

# Force column D to Text so numeric-looking price strings ("573.91",
# "0.997", ...) are written verbatim instead of becoming float cells.
$priceRange = $ws.Range("D$firstRow`:D$lastRow")
$priceRange.NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $firstRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
}

# Restore the default (unstyled) look for column D now that the text is in place.
$priceRange.Style = "Normal"

"Updated rows $firstRow to $lastRow"
